$wb = $excel.ActiveWorkbook

# Rename the second sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInc = $wb.Worksheets.Item("Include from FSIII")
$wsInc.Name = "Include #0"

# --- Metadata sheet updates ---
$wsMeta.Range("B3").Value = "1.2.0"
$wsMeta.Range("B7").Value = "false"
$wsMeta.Range("B8").Value = "2024-10-31T19:21:51+01:00"
$wsMeta.Range("B10").Value = "KL (http://www.kl.dk)"
$wsMeta.Range("B11").Value = ""

# --- Include sheet updates ---
$wsInc.Range("B2").Value = ""
$wsInc.Range("A3").Value = "d7ff926a-4955-478f-b300-0b0ec0785013"
$wsInc.Range("B3").Value = ""

$wsInc.Range("A5").Value = "System URI"
$wsInc.Range("B5").Value = "urn:oid:1.2.208.176.2.21"
$wsInc.Range("A4").Value = ""
$wsInc.Range("B4").Value = ""
